# Auto-generated script to update price/profit figures across all 8 leve-profit sheets.
# Values correspond to a scheduled market-data refresh (currentAveragePrice* / Leve*Price* / Leve*Profit* columns).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (55 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1134
$ws.Range("J17").Value = 1134
$ws.Range("L17").Value = 3402
$ws.Range("N17").Value = -3738
$ws.Range("H51").Value = 7964.85
$ws.Range("J51").Value = 3249.7856
$ws.Range("L51").Value = 3249.7856
$ws.Range("N51").Value = -4217.7856
$ws.Range("H53").Value = 481.84616
$ws.Range("I53").Value = 528.2727
$ws.Range("J53").Value = 447.8
$ws.Range("K53").Value = 528.2727
$ws.Range("L53").Value = 447.8
$ws.Range("M53").Value = 108.7273
$ws.Range("N53").Value = -1721.8
$ws.Range("H98").Value = 553.8889
$ws.Range("I98").Value = 581.6667
$ws.Range("K98").Value = 581.6667
$ws.Range("M98").Value = 916.3333
$ws.Range("H112").Value = 1923.909
$ws.Range("I112").Value = 1480
$ws.Range("J112").Value = 1968.3
$ws.Range("K112").Value = 4440
$ws.Range("L112").Value = 5904.9
$ws.Range("M112").Value = -3332
$ws.Range("N112").Value = -8120.9
$ws.Range("H122").Value = 553.8889
$ws.Range("I122").Value = 581.6667
$ws.Range("K122").Value = 1745.0001
$ws.Range("M122").Value = 704.9999
$ws.Range("H129").Value = 937.4337
$ws.Range("J129").Value = 973.9067
$ws.Range("L129").Value = 2921.7201
$ws.Range("N129").Value = -12921.7201
$ws.Range("H135").Value = 2253.375
$ws.Range("I135").Value = 660.625
$ws.Range("J135").Value = 3846.125
$ws.Range("K135").Value = 5945.625
$ws.Range("L135").Value = 34615.125
$ws.Range("M135").Value = -3410.625
$ws.Range("N135").Value = -39685.125
$ws.Range("H138").Value = 3025.4119
$ws.Range("I138").Value = 2948.2222
$ws.Range("J138").Value = 3112.25
$ws.Range("K138").Value = 8844.6666
$ws.Range("L138").Value = 9336.75
$ws.Range("M138").Value = -3704.6666
$ws.Range("N138").Value = -19616.75
$ws.Range("H141").Value = 2599.2144
$ws.Range("I141").Value = 2292.4167
$ws.Range("J141").Value = 4440
$ws.Range("K141").Value = 6877.250100000001
$ws.Range("L141").Value = 13320
$ws.Range("M141").Value = -1697.250100000001
$ws.Range("N141").Value = -23680

# ---- Sheet: ARM (49 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 10312.75
$ws.Range("I6").Value = 17234
$ws.Range("J6").Value = 6160
$ws.Range("K6").Value = 17234
$ws.Range("L6").Value = 6160
$ws.Range("M6").Value = -17061
$ws.Range("N6").Value = -6506
$ws.Range("H74").Value = 994.36584
$ws.Range("I74").Value = 993.34375
$ws.Range("J74").Value = 998
$ws.Range("K74").Value = 993.34375
$ws.Range("L74").Value = 998
$ws.Range("M74").Value = -119.34375
$ws.Range("N74").Value = -2746
$ws.Range("H77").Value = 994.36584
$ws.Range("I77").Value = 993.34375
$ws.Range("J77").Value = 998
$ws.Range("K77").Value = 4966.71875
$ws.Range("L77").Value = 4990
$ws.Range("M77").Value = -598.71875
$ws.Range("N77").Value = -13726
$ws.Range("H88").Value = 1185.2858
$ws.Range("I88").Value = 1179.6
$ws.Range("J88").Value = 1199.5
$ws.Range("K88").Value = 1179.6
$ws.Range("L88").Value = 1199.5
$ws.Range("M88").Value = -773.5999999999999
$ws.Range("N88").Value = -2011.5
$ws.Range("H91").Value = 1185.2858
$ws.Range("I91").Value = 1179.6
$ws.Range("J91").Value = 1199.5
$ws.Range("K91").Value = 1179.6
$ws.Range("L91").Value = 1199.5
$ws.Range("M91").Value = 224.4000000000001
$ws.Range("N91").Value = -4007.5
$ws.Range("H97").Value = 34289.133
$ws.Range("I97").Value = 50575
$ws.Range("J97").Value = 1717.4
$ws.Range("K97").Value = 50575
$ws.Range("L97").Value = 1717.4
$ws.Range("M97").Value = -50079
$ws.Range("N97").Value = -2709.4
$ws.Range("H132").Value = 1982.6316
$ws.Range("I132").Value = 1564.6666
$ws.Range("J132").Value = 3550
$ws.Range("K132").Value = 4693.9998
$ws.Range("L132").Value = 10650
$ws.Range("M132").Value = -2163.9998
$ws.Range("N132").Value = -15710

# ---- Sheet: BSM (21 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 88160.84
$ws.Range("I86").Value = 162083.72
$ws.Range("J86").Value = 1917.5
$ws.Range("K86").Value = 162083.72
$ws.Range("L86").Value = 1917.5
$ws.Range("M86").Value = -160960.72
$ws.Range("N86").Value = -4163.5
$ws.Range("H89").Value = 88160.84
$ws.Range("I89").Value = 162083.72
$ws.Range("J89").Value = 1917.5
$ws.Range("K89").Value = 810418.6
$ws.Range("L89").Value = 9587.5
$ws.Range("M89").Value = -804802.6
$ws.Range("N89").Value = -20819.5
$ws.Range("H94").Value = 559
$ws.Range("I94").Value = 488.4
$ws.Range("J94").Value = 710.2857
$ws.Range("K94").Value = 488.4
$ws.Range("L94").Value = 710.2857
$ws.Range("M94").Value = -37.39999999999998
$ws.Range("N94").Value = -1612.2857

# ---- Sheet: CRP (49 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 679
$ws.Range("I16").Value = 832.6667
$ws.Range("J16").Value = 525.3333
$ws.Range("K16").Value = 832.6667
$ws.Range("L16").Value = 525.3333
$ws.Range("M16").Value = -545.6667
$ws.Range("N16").Value = -1099.3333
$ws.Range("H31").Value = 2151.2793
$ws.Range("I31").Value = 1355.5151
$ws.Range("J31").Value = 2901.5715
$ws.Range("K31").Value = 1355.5151
$ws.Range("L31").Value = 2901.5715
$ws.Range("M31").Value = -1060.5151
$ws.Range("N31").Value = -3491.5715
$ws.Range("H34").Value = 2151.2793
$ws.Range("I34").Value = 1355.5151
$ws.Range("J34").Value = 2901.5715
$ws.Range("K34").Value = 1355.5151
$ws.Range("L34").Value = 2901.5715
$ws.Range("M34").Value = -1153.5151
$ws.Range("N34").Value = -3305.5715
$ws.Range("H62").Value = 2610
$ws.Range("I62").Value = 2490
$ws.Range("J62").Value = 2625
$ws.Range("K62").Value = 2490
$ws.Range("L62").Value = 2625
$ws.Range("M62").Value = -1866
$ws.Range("N62").Value = -3873
$ws.Range("H65").Value = 2610
$ws.Range("I65").Value = 2490
$ws.Range("J65").Value = 2625
$ws.Range("K65").Value = 12450
$ws.Range("L65").Value = 13125
$ws.Range("M65").Value = -9330
$ws.Range("N65").Value = -19365
$ws.Range("H113").Value = 679
$ws.Range("I113").Value = 832.6667
$ws.Range("J113").Value = 525.3333
$ws.Range("K113").Value = 832.6667
$ws.Range("L113").Value = 525.3333
$ws.Range("M113").Value = 1337.3333
$ws.Range("N113").Value = -4865.3333
$ws.Range("H132").Value = 5579.75
$ws.Range("I132").Value = 6316.4165
$ws.Range("J132").Value = 4474.75
$ws.Range("K132").Value = 18949.2495
$ws.Range("L132").Value = 13424.25
$ws.Range("M132").Value = -16419.2495
$ws.Range("N132").Value = -18484.25

# ---- Sheet: CUL (14 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3026.75
$ws.Range("I109").Value = 887.55554
$ws.Range("J109").Value = 4310.2666
$ws.Range("K109").Value = 2662.66662
$ws.Range("L109").Value = 12930.7998
$ws.Range("M109").Value = -1622.66662
$ws.Range("N109").Value = -15010.7998
$ws.Range("H131").Value = 821.89
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 827.1616
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2481.4848
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12561.4848

# ---- Sheet: GSM (8 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 91004740
$ws.Range("I80").Value = 166839800
$ws.Range("K80").Value = 166839800
$ws.Range("M80").Value = -166838802
$ws.Range("H83").Value = 91004740
$ws.Range("I83").Value = 166839800
$ws.Range("K83").Value = 834199000
$ws.Range("M83").Value = -834194008

# ---- Sheet: LTW (15 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 73578.36
$ws.Range("I40").Value = 334300
$ws.Range("J40").Value = 2472.4546
$ws.Range("K40").Value = 334300
$ws.Range("L40").Value = 2472.4546
$ws.Range("M40").Value = -334164
$ws.Range("N40").Value = -2744.4546
$ws.Range("H81").Value = 37160
$ws.Range("J81").Value = 37160
$ws.Range("L81").Value = 37160
$ws.Range("N81").Value = -39156
$ws.Range("H84").Value = 37160
$ws.Range("J84").Value = 37160
$ws.Range("L84").Value = 111480
$ws.Range("N84").Value = -121464

# ---- Sheet: WVR (33 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12822762
$ws.Range("J62").Value = 2700
$ws.Range("L62").Value = 2700
$ws.Range("N62").Value = -3948
$ws.Range("H65").Value = 12822762
$ws.Range("J65").Value = 2700
$ws.Range("L65").Value = 13500
$ws.Range("N65").Value = -19740
$ws.Range("H81").Value = 154364.39
$ws.Range("I81").Value = 125322.25
$ws.Range("J81").Value = 200831.8
$ws.Range("K81").Value = 250644.5
$ws.Range("L81").Value = 401663.6
$ws.Range("M81").Value = -249583.5
$ws.Range("N81").Value = -403785.6
$ws.Range("H84").Value = 154364.39
$ws.Range("I84").Value = 125322.25
$ws.Range("J84").Value = 200831.8
$ws.Range("K84").Value = 1253222.5
$ws.Range("L84").Value = 2008318
$ws.Range("M84").Value = -1247918.5
$ws.Range("N84").Value = -2018926
$ws.Range("H100").Value = 112567.11
$ws.Range("I100").Value = 144114.86
$ws.Range("J100").Value = 2150
$ws.Range("K100").Value = 288229.72
$ws.Range("L100").Value = 4300
$ws.Range("M100").Value = -287688.72
$ws.Range("N100").Value = -5382
$ws.Range("H122").Value = 2499.25
$ws.Range("I122").Value = 1923.5
$ws.Range("K122").Value = 5770.5
$ws.Range("M122").Value = -3320.5
